$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$row4 = @("f4 address", "f4 city", "f4 first", "f4 lasst ", "040404", "0404", "0404")
$row5 = @("add", "city", "firsy", "last", " ", "033", "333")

for ($i = 0; $i -lt $row4.Length; $i++) {
    $cell = $ws.Cells.Item(4, $i + 1)
    $cell.NumberFormat = "@"
    $cell.Value = $row4[$i]
    $cell.Style = "Normal"
}

for ($i = 0; $i -lt $row5.Length; $i++) {
    $cell = $ws.Cells.Item(5, $i + 1)
    $cell.NumberFormat = "@"
    $cell.Value = $row5[$i]
    $cell.Style = "Normal"
}
